$wb = $excel.ActiveWorkbook

# --- Locate the existing "ODI Batting" sheet ---
$odiSheet = $wb.Worksheets.Item("ODI Batting")

# --- Insert a new "Player Info" sheet before it ---
$playerInfo = $wb.Worksheets.Add($odiSheet)
$playerInfo.Name = "Player Info"

# --- Header row (bold, centered, bordered - same look as other sheets) ---
$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data row ---
# The leading apostrophe forces text (not numeric) storage for the numeric-
# looking ID; re-applying the "Normal" style afterwards drops the stray
# quote-prefix formatting so the cell keeps the workbook's plain default look.
$playerInfo.Cells.Item(2,1).Value = "'5659"
$playerInfo.Cells.Item(2,1).Style = "Normal"
$playerInfo.Cells.Item(2,2).Value = "Mohammad Naim Sheikh"
$playerInfo.Cells.Item(2,3).Value = "Left Handed"
$playerInfo.Cells.Item(2,4).Value = "Does Not Bowl | Unknown"

# --- Update the "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE,
#     and replace full scorecard URLs with just the numeric match code.
#     Re-fetch the sheet by name now that the insert has happened - the
#     handle captured before Worksheets.Add() can reseat onto the newly
#     inserted sheet instead of the original one. ---
$odiSheet = $wb.Worksheets.Item("ODI Batting")
$odiSheet.Cells.Item(1,4).Value = "MATCH_CODE"
$odiSheet.Cells.Item(2,4).Value = "'4420"
$odiSheet.Cells.Item(2,4).Style = "Normal"
$odiSheet.Cells.Item(3,4).Value = "'4465"
$odiSheet.Cells.Item(3,4).Style = "Normal"
